# Refresh market-price-derived columns (H:N) on each Leve-profit sheet
# with the latest pull from the pricing source. Values only; no
# structural changes.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1815.0834
$ws.Range("I6").Value = 176.2
$ws.Range("K6").Value = 528.5999999999999
$ws.Range("M6").Value = -416.5999999999999
$ws.Range("H51").Value = 20001780
$ws.Range("I51").Value = 1983.5
$ws.Range("J51").Value = 50001476
$ws.Range("K51").Value = 1983.5
$ws.Range("L51").Value = 50001476
$ws.Range("M51").Value = -1499.5
$ws.Range("N51").Value = -50002444
$ws.Range("H70").Value = 1720.6072
$ws.Range("I70").Value = 1210.579
$ws.Range("J70").Value = 2797.3333
$ws.Range("K70").Value = 3631.737
$ws.Range("L70").Value = 8391.999899999999
$ws.Range("M70").Value = -3361.737
$ws.Range("N70").Value = -8931.999899999999
$ws.Range("H73").Value = 1720.6072
$ws.Range("I73").Value = 1210.579
$ws.Range("J73").Value = 2797.3333
$ws.Range("K73").Value = 3631.737
$ws.Range("L73").Value = 8391.999899999999
$ws.Range("M73").Value = -2695.737
$ws.Range("N73").Value = -10263.9999
$ws.Range("H76").Value = 2711.35
$ws.Range("I76").Value = 2612.0557
$ws.Range("J76").Value = 3605
$ws.Range("K76").Value = 2612.0557
$ws.Range("L76").Value = 3605
$ws.Range("M76").Value = -2297.0557
$ws.Range("N76").Value = -4235
$ws.Range("H79").Value = 2711.35
$ws.Range("I79").Value = 2612.0557
$ws.Range("J79").Value = 3605
$ws.Range("K79").Value = 2612.0557
$ws.Range("L79").Value = 3605
$ws.Range("M79").Value = -1520.0557
$ws.Range("N79").Value = -5789
$ws.Range("H92").Value = 286.125
$ws.Range("I92").Value = 269.83334
$ws.Range("J92").Value = 335
$ws.Range("K92").Value = 269.83334
$ws.Range("L92").Value = 335
$ws.Range("M92").Value = 978.16666
$ws.Range("N92").Value = -2831
$ws.Range("H121").Value = 1334.5454
$ws.Range("J121").Value = 1622.5
$ws.Range("L121").Value = 4867.5
$ws.Range("N121").Value = -8361.5
$ws.Range("H125").Value = 1853.8823
$ws.Range("I125").Value = 1616.6666
$ws.Range("J125").Value = 1983.2727
$ws.Range("K125").Value = 14549.9994
$ws.Range("L125").Value = 17849.4543
$ws.Range("M125").Value = -12089.9994
$ws.Range("N125").Value = -22769.4543
$ws.Range("H141").Value = 643445.1
$ws.Range("I141").Value = 2078.5386
$ws.Range("K141").Value = 6235.6158
$ws.Range("M141").Value = -1055.6158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4269.88
$ws.Range("I32").Value = 4140.5
$ws.Range("J32").Value = 7375
$ws.Range("K32").Value = 4140.5
$ws.Range("L32").Value = 7375
$ws.Range("M32").Value = -3853.5
$ws.Range("N32").Value = -7949
$ws.Range("H35").Value = 7037
$ws.Range("I35").Value = 7037
$ws.Range("K35").Value = 7037
$ws.Range("M35").Value = -6631
$ws.Range("H45").Value = 1470.3334
$ws.Range("I45").Value = 1105.0385
$ws.Range("J45").Value = 2827.1428
$ws.Range("K45").Value = 1105.0385
$ws.Range("L45").Value = 2827.1428
$ws.Range("M45").Value = -728.0385000000001
$ws.Range("N45").Value = -3581.1428
$ws.Range("H132").Value = 2409.7932
$ws.Range("I132").Value = 2137.7234
$ws.Range("J132").Value = 3572.2727
$ws.Range("K132").Value = 6413.1702
$ws.Range("L132").Value = 10716.8181
$ws.Range("M132").Value = -3883.1702
$ws.Range("N132").Value = -15776.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 38778
$ws.Range("I86").Value = 3111.875
$ws.Range("J86").Value = 86332.836
$ws.Range("K86").Value = 3111.875
$ws.Range("L86").Value = 86332.836
$ws.Range("M86").Value = -1988.875
$ws.Range("N86").Value = -88578.836
$ws.Range("H89").Value = 38778
$ws.Range("I89").Value = 3111.875
$ws.Range("J89").Value = 86332.836
$ws.Range("K89").Value = 15559.375
$ws.Range("L89").Value = 431664.18
$ws.Range("M89").Value = -9943.375
$ws.Range("N89").Value = -442896.18
$ws.Range("H107").Value = 3354.2727
$ws.Range("I107").Value = 732.8333
$ws.Range("J107").Value = 6500
$ws.Range("K107").Value = 732.8333
$ws.Range("L107").Value = 6500
$ws.Range("M107").Value = 1187.1667
$ws.Range("N107").Value = -10340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2155
$ws.Range("J16").Value = 2726.2222
$ws.Range("L16").Value = 2726.2222
$ws.Range("N16").Value = -3300.2222
$ws.Range("H31").Value = 3678.6216
$ws.Range("I31").Value = 2590.9565
$ws.Range("J31").Value = 5465.5
$ws.Range("K31").Value = 2590.9565
$ws.Range("L31").Value = 5465.5
$ws.Range("M31").Value = -2295.9565
$ws.Range("N31").Value = -6055.5
$ws.Range("H34").Value = 3678.6216
$ws.Range("I34").Value = 2590.9565
$ws.Range("J34").Value = 5465.5
$ws.Range("K34").Value = 2590.9565
$ws.Range("L34").Value = 5465.5
$ws.Range("M34").Value = -2388.9565
$ws.Range("N34").Value = -5869.5
$ws.Range("H94").Value = 5292.8
$ws.Range("J94").Value = 4629.6
$ws.Range("L94").Value = 4629.6
$ws.Range("N94").Value = -5531.6
$ws.Range("H107").Value = 1268.9231
$ws.Range("I107").Value = 518.3
$ws.Range("K107").Value = 518.3
$ws.Range("M107").Value = 1401.7
$ws.Range("H113").Value = 2155
$ws.Range("J113").Value = 2726.2222
$ws.Range("L113").Value = 2726.2222
$ws.Range("N113").Value = -7066.2222
$ws.Range("H134").Value = 2486.4211
$ws.Range("I134").Value = 627.7778
$ws.Range("J134").Value = 4159.2
$ws.Range("K134").Value = 1883.3334
$ws.Range("L134").Value = 12477.6
$ws.Range("M134").Value = 651.6666
$ws.Range("N134").Value = -17547.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1180.7865
$ws.Range("J131").Value = 1006.875
$ws.Range("L131").Value = 3020.625
$ws.Range("N131").Value = -13100.625
$ws.Range("H140").Value = 1787.48
$ws.Range("I140").Value = 921.5294
$ws.Range("J140").Value = 3627.625
$ws.Range("K140").Value = 2764.5882
$ws.Range("L140").Value = 10882.875
$ws.Range("M140").Value = 2415.4118
$ws.Range("N140").Value = -21242.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716
$ws.Range("H113").Value = 2006.65
$ws.Range("I113").Value = 1669.4166
$ws.Range("J113").Value = 2512.5
$ws.Range("K113").Value = 1669.4166
$ws.Range("L113").Value = 2512.5
$ws.Range("M113").Value = 500.5834
$ws.Range("N113").Value = -6852.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2671.111
$ws.Range("I100").Value = 1500
$ws.Range("J100").Value = 2817.5
$ws.Range("K100").Value = 1500
$ws.Range("L100").Value = 2817.5
$ws.Range("M100").Value = -959
$ws.Range("N100").Value = -3899.5
$ws.Range("H132").Value = 4276.533
$ws.Range("I132").Value = 3024.6667
$ws.Range("K132").Value = 9074.000100000001
$ws.Range("M132").Value = -6544.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4409.383
$ws.Range("I132").Value = 1848.2368
$ws.Range("J132").Value = 15223.111
$ws.Range("K132").Value = 5544.7104
$ws.Range("L132").Value = 45669.333
$ws.Range("M132").Value = -3014.7104
$ws.Range("N132").Value = -50729.333
